# Auto-generated edit script applying numeric corrections to the leve-profit
# calculation columns (H:N) on each class sheet, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2326.9546
$ws.Range("I2").Value = 619.8889
$ws.Range("J2").Value = 3508.7693
$ws.Range("K2").Value = 619.8889
$ws.Range("L2").Value = 3508.7693
$ws.Range("M2").Value = -506.8889
$ws.Range("N2").Value = -3734.7693
$ws.Range("H33").Value = 272.4
$ws.Range("I33").Value = 177.5
$ws.Range("J33").Value = 335.66666
$ws.Range("K33").Value = 177.5
$ws.Range("L33").Value = 335.66666
$ws.Range("M33").Value = 51.5
$ws.Range("N33").Value = -793.66666
$ws.Range("H51").Value = 11187.333
$ws.Range("I51").Value = 5685.3335
$ws.Range("K51").Value = 5685.3335
$ws.Range("M51").Value = -5201.3335
$ws.Range("H62").Value = 7936.7144
$ws.Range("I62").Value = 7900
$ws.Range("J62").Value = 7942.8335
$ws.Range("K62").Value = 7900
$ws.Range("L62").Value = 7942.8335
$ws.Range("M62").Value = -7276
$ws.Range("N62").Value = -9190.833500000001
$ws.Range("H65").Value = 7936.7144
$ws.Range("I65").Value = 7900
$ws.Range("J65").Value = 7942.8335
$ws.Range("K65").Value = 39500
$ws.Range("L65").Value = 39714.1675
$ws.Range("M65").Value = -36380
$ws.Range("N65").Value = -45954.1675
$ws.Range("H70").Value = 8757.25
$ws.Range("I70").Value = 5773.75
$ws.Range("J70").Value = 14724.25
$ws.Range("K70").Value = 17321.25
$ws.Range("L70").Value = 44172.75
$ws.Range("M70").Value = -17051.25
$ws.Range("N70").Value = -44712.75
$ws.Range("H73").Value = 8757.25
$ws.Range("I73").Value = 5773.75
$ws.Range("J73").Value = 14724.25
$ws.Range("K73").Value = 17321.25
$ws.Range("L73").Value = 44172.75
$ws.Range("M73").Value = -16385.25
$ws.Range("N73").Value = -46044.75
$ws.Range("H74").Value = 138434.8
$ws.Range("I74").Value = 152982.67
$ws.Range("J74").Value = 7504
$ws.Range("K74").Value = 152982.67
$ws.Range("L74").Value = 7504
$ws.Range("M74").Value = -152046.67
$ws.Range("N74").Value = -9376
$ws.Range("H77").Value = 138434.8
$ws.Range("I77").Value = 152982.67
$ws.Range("J77").Value = 7504
$ws.Range("K77").Value = 764913.3500000001
$ws.Range("L77").Value = 37520
$ws.Range("M77").Value = -760233.3500000001
$ws.Range("N77").Value = -46880
$ws.Range("H92").Value = 387.66666
$ws.Range("J92").Value = 420.57144
$ws.Range("L92").Value = 420.57144
$ws.Range("N92").Value = -2916.57144
$ws.Range("H100").Value = 1796.6666
$ws.Range("I100").Value = 195
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 195
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = 346
$ws.Range("N100").Value = -6082
$ws.Range("H106").Value = 994.5
$ws.Range("I106").Value = 994.5
$ws.Range("K106").Value = 994.5
$ws.Range("M106").Value = -363.5
$ws.Range("H137").Value = 2021.625
$ws.Range("I137").Value = 2039
$ws.Range("K137").Value = 6117
$ws.Range("M137").Value = -3567

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 298.33334
$ws.Range("I25").Value = 298.33334
$ws.Range("K25").Value = 298.33334
$ws.Range("M25").Value = 103.66666
$ws.Range("H61").Value = 4300.773
$ws.Range("I61").Value = 1957.125
$ws.Range("K61").Value = 1957.125
$ws.Range("M61").Value = -1745.125
$ws.Range("H74").Value = 2272.0476
$ws.Range("I74").Value = 1788.8125
$ws.Range("J74").Value = 3818.4
$ws.Range("K74").Value = 1788.8125
$ws.Range("L74").Value = 3818.4
$ws.Range("M74").Value = -914.8125
$ws.Range("N74").Value = -5566.4
$ws.Range("H77").Value = 2272.0476
$ws.Range("I77").Value = 1788.8125
$ws.Range("J77").Value = 3818.4
$ws.Range("K77").Value = 8944.0625
$ws.Range("L77").Value = 19092
$ws.Range("M77").Value = -4576.0625
$ws.Range("N77").Value = -27828
$ws.Range("H102").Value = 1249.5
$ws.Range("I102").Value = 1249.5
$ws.Range("K102").Value = 1249.5
$ws.Range("M102").Value = 372.5
$ws.Range("H132").Value = 1785.7858
$ws.Range("I132").Value = 1785.7858
$ws.Range("K132").Value = 5357.357400000001
$ws.Range("M132").Value = -2827.357400000001
$ws.Range("H136").Value = 4300.773
$ws.Range("I136").Value = 1957.125
$ws.Range("K136").Value = 5871.375
$ws.Range("M136").Value = -3321.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H95").Value = 15874.667
$ws.Range("J95").Value = 15874.667
$ws.Range("L95").Value = 15874.667
$ws.Range("N95").Value = -21366.667
$ws.Range("H99").Value = 2621
$ws.Range("I99").Value = 2307.8333
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 2307.8333
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -809.8332999999998
$ws.Range("N99").Value = -7496
$ws.Range("H134").Value = 1602.4445
$ws.Range("I134").Value = 1602.4445
$ws.Range("K134").Value = 4807.333500000001
$ws.Range("M134").Value = -2272.333500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1953
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 1906
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 1906
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3154
$ws.Range("H65").Value = 1953
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 1906
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 9530
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -15770

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 42857340
$ws.Range("I23").Value = 100000080
$ws.Range("J23").Value = 286
$ws.Range("K23").Value = 300000240
$ws.Range("L23").Value = 858
$ws.Range("M23").Value = -300000005
$ws.Range("N23").Value = -1328
$ws.Range("H131").Value = 2897.1428
$ws.Range("J131").Value = 2897.1428
$ws.Range("L131").Value = 8691.428400000001
$ws.Range("N131").Value = -18771.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 25000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 25000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 25000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -25346
$ws.Range("H113").Value = 1143.875
$ws.Range("I113").Value = 1143.875
$ws.Range("K113").Value = 1143.875
$ws.Range("M113").Value = 1026.125
$ws.Range("H126").Value = 2999.5
$ws.Range("J126").Value = 3249
$ws.Range("L126").Value = 9747
$ws.Range("N126").Value = -14687
$ws.Range("H132").Value = 882.8461
$ws.Range("I132").Value = 873.0833
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2619.2499
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -89.2498999999998
$ws.Range("N132").Value = -8060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 824.5454999999999
$ws.Range("I16").Value = 422.22223
$ws.Range("J16").Value = 2635
$ws.Range("K16").Value = 422.22223
$ws.Range("L16").Value = 2635
$ws.Range("M16").Value = -252.22223
$ws.Range("N16").Value = -2975
$ws.Range("H68").Value = 2663.375
$ws.Range("I68").Value = 1832.6666
$ws.Range("K68").Value = 1832.6666
$ws.Range("M68").Value = -1083.6666
$ws.Range("H71").Value = 2663.375
$ws.Range("I71").Value = 1832.6666
$ws.Range("K71").Value = 9163.333000000001
$ws.Range("M71").Value = -5419.333000000001
$ws.Range("H82").Value = 1439.25
$ws.Range("I82").Value = 1568.5
$ws.Range("K82").Value = 1568.5
$ws.Range("M82").Value = -1207.5
$ws.Range("H85").Value = 1439.25
$ws.Range("I85").Value = 1568.5
$ws.Range("K85").Value = 1568.5
$ws.Range("M85").Value = -320.5
$ws.Range("H122").Value = 7335.607
$ws.Range("I122").Value = 7527.8887
$ws.Range("J122").Value = 6989.5
$ws.Range("K122").Value = 22583.6661
$ws.Range("L122").Value = 20968.5
$ws.Range("M122").Value = -20133.6661
$ws.Range("N122").Value = -25868.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 779.6
$ws.Range("I132").Value = 774.5
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 2323.5
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = 206.5
$ws.Range("N132").Value = -7460
